# The document contains "simple field" (<w:fldSimple>) usages that were
# authored with the wrong Word XML syntax (m2doc issue: fields must be
# the expanded/"complex" field form - begin/instrText/separate/end -
# otherwise user edits to the field result are lost on save/regeneration).
#
# Word itself performs exactly this normalization the first time such a
# document is opened and saved again: every <w:fldSimple w:instr="..."/>
# is unwound into the four-run complex-field sequence:
#   <w:r><w:fldChar w:fldCharType="begin"/></w:r>
#   <w:r><w:instrText>INSTR</w:instrText></w:r>
#   <w:r><w:fldChar w:fldCharType="separate"/></w:r>
#   <w:r><w:fldChar w:fldCharType="end"/></w:r>
#
# Reproduce that normalization for every simple field in the document by
# replacing each field's paragraph content with the expanded form via
# Range.InsertXML, while preserving any bookmarks already anchored in
# that paragraph (e.g. the automatic "_GoBack" bookmark), which would
# otherwise be lost since InsertXML replaces the whole paragraph.

$d = $word.ActiveDocument

$xmlTemplate = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>__PREFIX__<w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText>__INSTR__</w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:p></w:body></w:document>
</pkg:xmlData></pkg:part></pkg:package>
'@

# Walk the fields back-to-front so earlier InsertXML calls (which change
# run/paragraph content) never invalidate the indices of fields still to
# be processed.
for ($i = $d.Fields.Count; $i -ge 1; $i--) {
    $field = $d.Fields.Item($i)
    $instr = $field.Code.Text.Trim()
    # Escape for safe embedding inside the <w:instrText> element below.
    $instr = $instr.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $p = $field.Code.Paragraphs(1)

    # Pull this paragraph's own raw OOXML so any markup that precedes the
    # field (bookmarks, in this document) can be carried over verbatim -
    # the Bookmarks collection does not surface hidden bookmarks such as
    # "_GoBack", so read the XML directly instead.
    $raw = $p.Range.WordOpenXML
    $prefix = ""
    if ($raw -match '(?s)<w:p[ >][^>]*>(.*?)<w:fldSimple\b') {
        $prefix = $Matches[1]
    }

    $xml = $xmlTemplate.Replace("__PREFIX__", $prefix).Replace("__INSTR__", $instr)
    [void]$p.Range.InsertXML($xml)
}
